# Updated cryptos list with GitHub Actions
# Refreshes Price/Volume(1h) figures scraped for the crypto table and fixes
# the TrustWalletToken / Algorand row ordering (rows 39-40).
# D-column cells hold text (e.g. "1.002", "24.542.09") rather than numbers,
# so NumberFormat is forced to "@" before each Value assignment to keep
# Excel from silently re-typing them as numeric and dropping trailing
# zeros / separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.594.52"
$ws.Range("E2").Value = "  -0.24%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.58"
$ws.Range("E3").Value = "  -2.54%  "

# --- Row 4: TetherUSD ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.04"
$ws.Range("E5").Value = "  -0.23%  "

# --- Row 6: USDC ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  -0.04%  "

# --- Row 7: XRP ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3620"
$ws.Range("E7").Value = "  -2.78%  "

# --- Row 8: OKB ---
$ws.Range("E8").Value = "  -2.58%  "

# --- Row 9: Cardano ---
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3259"
$ws.Range("E9").Value = "  -5.18%  "

# --- Row 10: Polygon ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  -4.56%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06966"
$ws.Range("E11").Value = "  -6.34%  "

# --- Row 12: BinanceUSD ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9972"
$ws.Range("E12").Value = "  -0.08%  "

# --- Row 13: Polkadot ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.907"
$ws.Range("E13").Value = "  -4.88%  "

# --- Row 14: Solana ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.43"
$ws.Range("E14").Value = "  -6.76%  "

# --- Row 15: WrappedEther ---
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.659.29"
$ws.Range("E15").Value = "  -2.78%  "

# --- Row 16: Chainlink ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.560"
$ws.Range("E16").Value = "  -5.19%  "

# --- Row 17: ShibaInu ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001049"

# --- Row 18: TRON ---
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06537"
$ws.Range("E18").Value = "  -2.29%  "

# --- Row 19: Dai ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9964"
$ws.Range("E19").Value = "  -0.15%  "

# --- Row 20: Litecoin ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.65"
$ws.Range("E20").Value = "  -7.87%  "

# --- Row 21: Uniswap ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.921"
$ws.Range("E21").Value = "  -6.54%  "

# --- Row 22: Avalanche ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.71"
$ws.Range("E22").Value = "  -7.78%  "

# --- Row 23: Cosmos ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.60"
$ws.Range("E23").Value = "  -4.04%  "

# --- Row 24: WrappedBTC ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.550.13"
$ws.Range("E24").Value = "  -0.29%  "

# --- Row 25: Toncoin ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.465"
$ws.Range("E25").Value = "  +2.20%  "

# --- Row 26: LidoDAOToken ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.310"
$ws.Range("E26").Value = "  -16.32%  "

# --- Row 27: Monero ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.75"
$ws.Range("E27").Value = "  -1.77%  "

# --- Row 28: EthereumClassic ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.50"

# --- Row 29: WrappedliquidstakedEther2.0 ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.845.17"
$ws.Range("E29").Value = "  -2.63%  "

# --- Row 30: ImmutableX ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.190"
$ws.Range("E30").Value = "  +1.64%  "

# --- Row 31: BitcoinCash ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.10"
$ws.Range("E31").Value = "  -5.00%  "

# --- Row 32: HuobiToken ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.059"
$ws.Range("E32").Value = "  -3.28%  "

# --- Row 33: Filecoin ---
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.644"
$ws.Range("E33").Value = "  -16.02%  "

# --- Row 34: WEMIXTOKEN ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.696"
$ws.Range("E34").Value = "  -4.03%  "

# --- Row 35: Stellar ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08364"
$ws.Range("E35").Value = "  -4.69%  "

# --- Row 36: Aptos ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.41"
$ws.Range("E36").Value = "  -9.14%  "

# --- Row 37: InternetComputer(DFINITY) ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.202"
$ws.Range("E37").Value = "  -5.35%  "

# --- Row 38: Hedera ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06067"
$ws.Range("E38").Value = "  -6.83%  "

# --- Rows 39-40: fix TrustWalletToken / Algorand order + refresh values ---
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2060"
$ws.Range("E39").Value = "  -6.99%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.207"
$ws.Range("E40").Value = "  -5.05%  "

# --- Row 41: FraxShare ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.234"

# --- Row 42: VeChain ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02192"
$ws.Range("E42").Value = "  -7.11%  "

# --- Row 43: Frax ---
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9975"

# --- Row 44: TheSandbox ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5897"
$ws.Range("E44").Value = "  -7.62%  "

# --- Row 45: PancakeSwap ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.737"
$ws.Range("E45").Value = "  -1.63%  "

# --- Row 46: EnergySwap ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.67"
$ws.Range("E46").Value = "  -8.37%  "

# --- Row 47: Decentraland ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5598"
$ws.Range("E47").Value = "  -7.65%  "

# --- Row 48: Quant ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.38"
$ws.Range("E48").Value = "  -4.95%  "

# --- Row 49: NEARProtocol ---
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.942"
$ws.Range("E49").Value = "  -7.84%  "

# --- Row 50: Cronos ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06945"

# --- Row 51: Aave ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.28"
$ws.Range("E51").Value = "  -5.77%  "

